$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (column C) and montant_total (column D) values for the
# rows affected by the 2020-05-22 data refresh.

$ws.Cells.Item(2, 3).Value = 184639
$ws.Cells.Item(2, 4).Value = 238151732
$ws.Cells.Item(4, 3).Value = 184
$ws.Cells.Item(4, 4).Value = 266539
$ws.Cells.Item(6, 3).Value = 496
$ws.Cells.Item(6, 4).Value = 737537
$ws.Cells.Item(8, 3).Value = 76163
$ws.Cells.Item(8, 4).Value = 112389534
$ws.Cells.Item(9, 3).Value = 92
$ws.Cells.Item(9, 4).Value = 137790
$ws.Cells.Item(10, 3).Value = 35910
$ws.Cells.Item(10, 4).Value = 52144147
$ws.Cells.Item(13, 3).Value = 2054
$ws.Cells.Item(13, 4).Value = 2926526
$ws.Cells.Item(16, 3).Value = 2355
$ws.Cells.Item(16, 4).Value = 3319813
$ws.Cells.Item(17, 3).Value = 46861
$ws.Cells.Item(17, 4).Value = 59639843
$ws.Cells.Item(22, 3).Value = 21548
$ws.Cells.Item(22, 4).Value = 31758151
$ws.Cells.Item(24, 3).Value = 7124
$ws.Cells.Item(24, 4).Value = 10338308
$ws.Cells.Item(26, 3).Value = 844
$ws.Cells.Item(26, 4).Value = 1187491
$ws.Cells.Item(28, 3).Value = 713
$ws.Cells.Item(28, 4).Value = 1011983
$ws.Cells.Item(29, 3).Value = 59785
$ws.Cells.Item(29, 4).Value = 76794305
$ws.Cells.Item(32, 3).Value = 506
$ws.Cells.Item(32, 4).Value = 751471
$ws.Cells.Item(34, 3).Value = 29637
$ws.Cells.Item(34, 4).Value = 43736700
$ws.Cells.Item(36, 3).Value = 5641
$ws.Cells.Item(36, 4).Value = 8131907
$ws.Cells.Item(38, 3).Value = 753
$ws.Cells.Item(38, 4).Value = 1050945
$ws.Cells.Item(39, 3).Value = 787
$ws.Cells.Item(39, 4).Value = 1108671
$ws.Cells.Item(40, 3).Value = 41307
$ws.Cells.Item(40, 4).Value = 52487377
$ws.Cells.Item(45, 3).Value = 18563
$ws.Cells.Item(45, 4).Value = 27375187
$ws.Cells.Item(47, 3).Value = 6785
$ws.Cells.Item(47, 4).Value = 9849141
$ws.Cells.Item(50, 3).Value = 510
$ws.Cells.Item(50, 4).Value = 723768
$ws.Cells.Item(51, 3).Value = 11514
$ws.Cells.Item(51, 4).Value = 15085946
$ws.Cells.Item(55, 3).Value = 4247
$ws.Cells.Item(55, 4).Value = 6235424
$ws.Cells.Item(56, 3).Value = 2899
$ws.Cells.Item(56, 4).Value = 4225547
$ws.Cells.Item(59, 3).Value = 83903
$ws.Cells.Item(59, 4).Value = 106100949
$ws.Cells.Item(65, 3).Value = 41677
$ws.Cells.Item(65, 4).Value = 61547988
$ws.Cells.Item(68, 3).Value = 18463
$ws.Cells.Item(68, 4).Value = 26847621
$ws.Cells.Item(70, 3).Value = 1351
$ws.Cells.Item(70, 4).Value = 1946279
$ws.Cells.Item(71, 3).Value = 1034
$ws.Cells.Item(71, 4).Value = 1466852
$ws.Cells.Item(72, 3).Value = 13636
$ws.Cells.Item(72, 4).Value = 18487120
$ws.Cells.Item(76, 3).Value = 4266
$ws.Cells.Item(76, 4).Value = 6291862
$ws.Cells.Item(78, 3).Value = 3354
$ws.Cells.Item(78, 4).Value = 4844868
$ws.Cells.Item(80, 3).Value = 196
$ws.Cells.Item(80, 4).Value = 281311
$ws.Cells.Item(82, 3).Value = 2884
$ws.Cells.Item(82, 4).Value = 3970010
$ws.Cells.Item(83, 3).Value = 853
$ws.Cells.Item(83, 4).Value = 1253754
$ws.Cells.Item(85, 3).Value = 1034
$ws.Cells.Item(85, 4).Value = 1503016
$ws.Cells.Item(88, 3).Value = 85964
$ws.Cells.Item(88, 4).Value = 107396946
$ws.Cells.Item(92, 3).Value = 531
$ws.Cells.Item(92, 4).Value = 787277
$ws.Cells.Item(94, 3).Value = 35420
$ws.Cells.Item(94, 4).Value = 52209709
$ws.Cells.Item(96, 3).Value = 16664
$ws.Cells.Item(96, 4).Value = 24235235
$ws.Cells.Item(97, 3).Value = 665
$ws.Cells.Item(97, 4).Value = 906842
$ws.Cells.Item(98, 3).Value = 836
$ws.Cells.Item(98, 4).Value = 1176491
$ws.Cells.Item(100, 3).Value = 24638
$ws.Cells.Item(100, 4).Value = 33256848
$ws.Cells.Item(105, 3).Value = 9009
$ws.Cells.Item(105, 4).Value = 13284682
$ws.Cells.Item(106, 3).Value = 2300
$ws.Cells.Item(106, 4).Value = 3333482
$ws.Cells.Item(108, 3).Value = 206
$ws.Cells.Item(108, 4).Value = 297215
$ws.Cells.Item(109, 3).Value = 155
$ws.Cells.Item(109, 4).Value = 216566
$ws.Cells.Item(110, 3).Value = 9016
$ws.Cells.Item(110, 4).Value = 12014981
$ws.Cells.Item(112, 3).Value = 4147
$ws.Cells.Item(112, 4).Value = 6069347
$ws.Cells.Item(113, 3).Value = 2588
$ws.Cells.Item(113, 4).Value = 3743585
$ws.Cells.Item(116, 3).Value = 3036
$ws.Cells.Item(116, 4).Value = 4300309
$ws.Cells.Item(117, 3).Value = 767
$ws.Cells.Item(117, 4).Value = 1136943
$ws.Cells.Item(121, 3).Value = 54001
$ws.Cells.Item(121, 4).Value = 68611480
$ws.Cells.Item(125, 3).Value = 361
$ws.Cells.Item(125, 4).Value = 532262
$ws.Cells.Item(127, 3).Value = 22844
$ws.Cells.Item(127, 4).Value = 33687139
$ws.Cells.Item(129, 3).Value = 8046
$ws.Cells.Item(129, 4).Value = 11681847
$ws.Cells.Item(131, 3).Value = 691
$ws.Cells.Item(131, 4).Value = 973931
$ws.Cells.Item(132, 3).Value = 667
$ws.Cells.Item(132, 4).Value = 952140
$ws.Cells.Item(133, 3).Value = 143252
$ws.Cells.Item(133, 4).Value = 180513755
$ws.Cells.Item(138, 3).Value = 516
$ws.Cells.Item(138, 4).Value = 759059
$ws.Cells.Item(140, 3).Value = 57854
$ws.Cells.Item(140, 4).Value = 85192235
$ws.Cells.Item(143, 3).Value = 20549
$ws.Cells.Item(143, 4).Value = 29718867
$ws.Cells.Item(146, 3).Value = 2800
$ws.Cells.Item(146, 4).Value = 4004977
$ws.Cells.Item(148, 3).Value = 1901
$ws.Cells.Item(148, 4).Value = 2642617
$ws.Cells.Item(150, 3).Value = 153939
$ws.Cells.Item(150, 4).Value = 192368073
$ws.Cells.Item(157, 3).Value = 61678
$ws.Cells.Item(157, 4).Value = 90725037
$ws.Cells.Item(160, 3).Value = 31307
$ws.Cells.Item(160, 4).Value = 45463696
$ws.Cells.Item(163, 3).Value = 2518
$ws.Cells.Item(163, 4).Value = 3552328
$ws.Cells.Item(166, 3).Value = 2196
$ws.Cells.Item(166, 4).Value = 3063300
$ws.Cells.Item(169, 3).Value = 63720
$ws.Cells.Item(169, 4).Value = 81182093
$ws.Cells.Item(175, 3).Value = 33263
$ws.Cells.Item(175, 4).Value = 49029406
$ws.Cells.Item(177, 3).Value = 7515
$ws.Cells.Item(177, 4).Value = 10843394
$ws.Cells.Item(179, 3).Value = 1136
$ws.Cells.Item(179, 4).Value = 1636883
$ws.Cells.Item(181, 3).Value = 1017
$ws.Cells.Item(181, 4).Value = 1423151
$ws.Cells.Item(182, 3).Value = 151726
$ws.Cells.Item(182, 4).Value = 193113869
$ws.Cells.Item(184, 3).Value = 146
$ws.Cells.Item(184, 4).Value = 214399
$ws.Cells.Item(190, 3).Value = 61972
$ws.Cells.Item(190, 4).Value = 91327866
$ws.Cells.Item(192, 3).Value = 39317
$ws.Cells.Item(192, 4).Value = 57185302
$ws.Cells.Item(194, 3).Value = 1328
$ws.Cells.Item(194, 4).Value = 1891321
$ws.Cells.Item(196, 3).Value = 1804
$ws.Cells.Item(196, 4).Value = 2530519
$ws.Cells.Item(197, 3).Value = 237205
$ws.Cells.Item(197, 4).Value = 313038457
$ws.Cells.Item(202, 3).Value = 742
$ws.Cells.Item(202, 4).Value = 1105218
$ws.Cells.Item(204, 3).Value = 124156
$ws.Cells.Item(204, 4).Value = 183419377
$ws.Cells.Item(205, 3).Value = 216
$ws.Cells.Item(205, 4).Value = 321880
$ws.Cells.Item(207, 3).Value = 98794
$ws.Cells.Item(207, 4).Value = 144006525
$ws.Cells.Item(209, 3).Value = 1245
$ws.Cells.Item(209, 4).Value = 1756599
$ws.Cells.Item(211, 3).Value = 2297
$ws.Cells.Item(211, 4).Value = 3257922
